$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1055892638440206, 4.193923464361182, 1.237120070547445),
    @(-12.19494866425143, 0.3453140144921644, 11.02819687465457),
    @(3.551830255760532, -8.829065172880417, 8.992697445851405),
    @(4.123927570738847, -6.680597377273238, 5.334880648918858),
    @(4.32835804291493, 5.768243861648356, 1.005123424080155),
    @(-0.7777560342033203, 4.314871306689256, 2.766960350972315),
    @(-5.723839211014095, 1.048127250851321, 5.411107558124477),
    @(-2.712708458585516, 9.354416433370432, 0.1980249521866684),
    @(1.900770921752156, -1.904725146743598, -0.06296492072769189),
    @(6.883454844636765, -7.698698745583567, -1.359904374716438)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
